$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged), update values
$ws.Range("B3").Value = 0.9957218342903249
$ws.Range("C3").Value = 0.9956375856168691
$ws.Range("D3").Value = 0.9760214675634495

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9923110495414994
$ws.Range("C4").Value = 0.9921908449967436
$ws.Range("D4").Value = 0.9731564244732521

# Row 5: AdaBoostRegressor -> MLPRegressor, update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9989430829149922
$ws.Range("C5").Value = 0.9986733503275795
$ws.Range("D5").Value = 0.9982627847602258
